$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Update the "last updated" timestamp string (row 1, col A)
$ws.Range("A1").Value = "Datos actualizados a 23 de Abril de 2020 a las 08:52"

# Row 25: Israel
$ws.Range("A25").Value = "Israel"
$ws.Range("B25").Value = 14592
$ws.Range("C25").Value = 94
$ws.Range("D25").Value = 5334
$ws.Range("E25").Value = 9067
$ws.Range("F25").Value = 136
$ws.Range("G25").Value = 2
$ws.Range("H25").Value = 191

# Row 41: Ucrania
$ws.Range("A41").Value = "Ucrania"
$ws.Range("B41").Value = 7170
$ws.Range("C41").Value = 578
$ws.Range("D41").Value = 504
$ws.Range("E41").Value = 6479
$ws.Range("F41").Value = 45
$ws.Range("G41").Value = 13
$ws.Range("H41").Value = 187

# Row 42: Catar
$ws.Range("A42").Value = "Catar"
$ws.Range("B42").Value = 7141
$ws.Range("C42").Value = 0
$ws.Range("D42").Value = 689
$ws.Range("E42").Value = 6442
$ws.Range("F42").Value = 37
$ws.Range("G42").Value = 0
$ws.Range("H42").Value = 10

# Row 43: Chequia
$ws.Range("A43").Value = "Chequia"
$ws.Range("B43").Value = 7136
$ws.Range("C43").Value = 4
$ws.Range("D43").Value = 2002
$ws.Range("E43").Value = 4924
$ws.Range("F43").Value = 76
$ws.Range("G43").Value = 2
$ws.Range("H43").Value = 210

# Row 44: Serbia
$ws.Range("A44").Value = "Serbia"
$ws.Range("B44").Value = 7114
$ws.Range("C44").Value = 0
$ws.Range("D44").Value = 1025
$ws.Range("E44").Value = 5955
$ws.Range("F44").Value = 101
$ws.Range("G44").Value = 0
$ws.Range("H44").Value = 134

# Row 45: Filipinas
$ws.Range("A45").Value = "Filipinas"
$ws.Range("B45").Value = 6710
$ws.Range("C45").Value = 0
$ws.Range("D45").Value = 693
$ws.Range("E45").Value = 5571
$ws.Range("F45").Value = 1
$ws.Range("G45").Value = 0
$ws.Range("H45").Value = 446

# Row 46: Australia
$ws.Range("A46").Value = "Australia"
$ws.Range("B46").Value = 6660
$ws.Range("C46").Value = 11
$ws.Range("D46").Value = 5041
$ws.Range("E46").Value = 1544
$ws.Range("F46").Value = 46
$ws.Range("G46").Value = 1
$ws.Range("H46").Value = 75

# Row 110: Georgia
$ws.Range("A110").Value = "Georgia"
$ws.Range("B110").Value = 420
$ws.Range("C110").Value = 4
$ws.Range("D110").Value = 107
$ws.Range("E110").Value = 308
$ws.Range("F110").Value = 6
$ws.Range("G110").Value = 0
$ws.Range("H110").Value = 5

# Row 126: El Salvador
$ws.Range("A126").Value = "El Salvador"
$ws.Range("B126").Value = 250
$ws.Range("C126").Value = 13
$ws.Range("D126").Value = 67
$ws.Range("E126").Value = 175
$ws.Range("F126").Value = 2
$ws.Range("G126").Value = 1
$ws.Range("H126").Value = 8

# Row 128: Islas Feroe
$ws.Range("A128").Value = "Islas Feroe"
$ws.Range("B128").Value = 187
$ws.Range("C128").Value = 2
$ws.Range("D128").Value = 178
$ws.Range("E128").Value = 9
$ws.Range("F128").Value = 0
$ws.Range("G128").Value = 0
$ws.Range("H128").Value = 0

# Row 129: Congo
$ws.Range("A129").Value = "Congo"
$ws.Range("B129").Value = 186
$ws.Range("C129").Value = 0
$ws.Range("D129").Value = 16
$ws.Range("E129").Value = 164
$ws.Range("F129").Value = 0
$ws.Range("G129").Value = 0
$ws.Range("H129").Value = 6

# Row 136: Gibraltar
$ws.Range("A136").Value = "Gibraltar"
$ws.Range("B136").Value = 132
$ws.Range("C136").Value = 0
$ws.Range("D136").Value = 127
$ws.Range("E136").Value = 5
$ws.Range("F136").Value = 0
$ws.Range("G136").Value = 0
$ws.Range("H136").Value = 0
